# Update the row of "newly added iAuthor" test-candidate details.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# Client Id
$ws.Range("A2").Value = "wqgxm552"
# Candidate ID
$ws.Range("B2").Value = 23103082
# User Name
$ws.Range("C2").Value = "kgmhkvc53"
# Exam Password
$ws.Range("D2").Value = "m`$5ZG9%b"
# First Name
$ws.Range("F2").Value = "YSeiLfel"
# Last Name
$ws.Range("G2").Value = "iHCK"
